$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Senast uppdaterad" (C) column date for rows 2-27 from
# 2023-09-21 (45190) to 2023-09-23 (45192).
$ws.Range("C2:C27").Value = 45192
